# Runtime optimizacija pakeista i O3 ir atlikti smulkus programos papildymai
# Updates benchmark timing data (strategy 1 / 2 / 3 tables) and refreshes the
# active cell selection/view to the "3 strategija" header merge block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 27
$ws.Range("D27").Value = 0.0096252500000000001
$ws.Range("E27").Value = 0.0065814999999999997
$ws.Range("F27").Value = 0.052938199999999998
$ws.Range("G27").Value = 0.0293637
$ws.Range("H27").Value = 0.258162
$ws.Range("I27").Value = 0.21557499999999999
$ws.Range("J27").Value = 2.1775199999999999
$ws.Range("K27").Value = 2.1295500000000001
$ws.Range("L27").Value = 23.725100000000001
$ws.Range("M27").Value = 26.2197
# row 28
$ws.Range("D28").Value = 0.00041733299999999998
$ws.Range("E28").Value = 0.00020066600000000001
$ws.Range("F28").Value = 0.0012636699999999999
$ws.Range("G28").Value = 0.00113342
$ws.Range("H28").Value = 0.0071027499999999997
$ws.Range("I28").Value = 0.0209792
$ws.Range("J28").Value = 0.099339899999999995
$ws.Range("K28").Value = 0.65255799999999997
$ws.Range("L28").Value = 0.905532
$ws.Range("M28").Value = 12.5252
# row 29
$ws.Range("D29").Value = 0.00027845799999999998
$ws.Range("E29").Value = 0.00035633399999999998
$ws.Range("F29").Value = 0.0015719600000000001
$ws.Range("G29").Value = 0.00163142
$ws.Range("H29").Value = 0.015906300000000002
$ws.Range("I29").Value = 0.029835
$ws.Range("J29").Value = 0.240173
$ws.Range("K29").Value = 0.55680200000000002
$ws.Range("L29").Value = 5.5067700000000004
$ws.Range("M29").Value = 65.600800000000007
# row 30
$ws.Range("D30").Value = 0.0010889599999999999
$ws.Range("E30").Value = 0.000774458
$ws.Range("F30").Value = 0.0027775399999999998
$ws.Range("G30").Value = 0.0018125800000000001
$ws.Range("H30").Value = 0.015661499999999998
$ws.Range("I30").Value = 0.022657299999999998
$ws.Range("J30").Value = 0.188666
$ws.Range("K30").Value = 0.20315900000000001
$ws.Range("L30").Value = 2.0125199999999999
$ws.Range("M30").Value = 2.3549199999999999
# row 31
$ws.Range("D31").Value = 0.00068216700000000002
$ws.Range("E31").Value = 0.00061625000000000004
$ws.Range("F31").Value = 0.0021565
$ws.Range("G31").Value = 0.0020855800000000001
$ws.Range("H31").Value = 0.0151857
$ws.Range("I31").Value = 0.016297699999999998
$ws.Range("J31").Value = 0.13837099999999999
$ws.Range("K31").Value = 0.195767
$ws.Range("L31").Value = 2.8715600000000001
$ws.Range("M31").Value = 2.3601700000000001
# row 32
$ws.Range("D32").Value = 0.012154999999999999
$ws.Range("E32").Value = 0.0085753700000000006
$ws.Range("F32").Value = 0.060753399999999999
$ws.Range("G32").Value = 0.036058600000000003
$ws.Range("H32").Value = 0.31212499999999999
$ws.Range("I32").Value = 0.30540899999999999
$ws.Range("J32").Value = 2.84416
$ws.Range("K32").Value = 3.7379199999999999
$ws.Range("L32").Value = 35.021599999999999
$ws.Range("M32").Value = 109.06100000000001
# row 40
$ws.Range("D40").Value = 0.0091308299999999995
$ws.Range("E40").Value = 0.0025149199999999999
$ws.Range("F40").Value = 0.039048199999999998
$ws.Range("G40").Value = 0.0213361
$ws.Range("H40").Value = 0.25558599999999998
$ws.Range("I40").Value = 0.29861799999999999
$ws.Range("K40").Value = 2.1996000000000002
$ws.Range("M40").Value = 21.8919
# row 41
$ws.Range("D41").Value = 0.00046029099999999999
$ws.Range("E41").Value = 0.00006084
$ws.Range("F41").Value = 0.0013691199999999999
$ws.Range("G41").Value = 0.00087495800000000001
$ws.Range("H41").Value = 0.0086944599999999993
$ws.Range("I41").Value = 0.0215617
$ws.Range("K41").Value = 0.63495800000000002
$ws.Range("M41").Value = 14.884600000000001
# row 42
$ws.Range("D42").Value = 0.00178283
$ws.Range("E42").Value = 0.000042125000000000001
$ws.Range("F42").Value = 0.074233599999999997
$ws.Range("G42").Value = 0.00066045899999999996
$ws.Range("H42").Value = 5.7900900000000002
$ws.Range("I42").Value = 0.018127899999999999
$ws.Range("K42").Value = 0.207813
$ws.Range("M42").Value = 9.4475300000000004
# row 43
$ws.Range("D43").Value = 0.037849899999999999
$ws.Range("E43").Value = 0.00024458299999999999
$ws.Range("F43").Value = 0.0021755400000000001
$ws.Range("G43").Value = 0.0019497500000000001
$ws.Range("H43").Value = 0.016341499999999998
$ws.Range("I43").Value = 0.013631000000000001
$ws.Range("K43").Value = 0.15093100000000001
$ws.Range("M43").Value = 9.9357399999999991
# row 44
$ws.Range("D44").Value = 0.032362200000000001
$ws.Range("E44").Value = 0.00025133300000000001
$ws.Range("F44").Value = 0.00160133
$ws.Range("G44").Value = 0.0035526199999999998
$ws.Range("H44").Value = 0.0151515
$ws.Range("I44").Value = 0.021311799999999999
$ws.Range("K44").Value = 0.23924100000000001
$ws.Range("M44").Value = 14.861499999999999
# row 45
$ws.Range("D45").Value = 0.081677
$ws.Range("E45").Value = 0.0031289999999999998
$ws.Range("F45").Value = 0.118482
$ws.Range("G45").Value = 0.0284026
$ws.Range("H45").Value = 6.0860300000000001
$ws.Range("I45").Value = 0.37331199999999998
$ws.Range("K45").Value = 3.4326300000000001
$ws.Range("M45").Value = 70.997600000000006
# row 53
$ws.Range("D53").Value = 0.0092640399999999994
$ws.Range("E53").Value = 0.0065449999999999996
$ws.Range("F53").Value = 0.052044800000000002
$ws.Range("G53").Value = 0.027611
$ws.Range("H53").Value = 0.25426199999999999
$ws.Range("I53").Value = 0.34542800000000001
$ws.Range("J53").Value = 2.1598700000000002
$ws.Range("K53").Value = 2.1947899999999998
$ws.Range("L53").Value = 24.0488
$ws.Range("M53").Value = 22.1267
# row 54
$ws.Range("D54").Value = 0.00043108399999999998
$ws.Range("E54").Value = 0.00015129100000000001
$ws.Range("F54").Value = 0.0012318299999999999
$ws.Range("G54").Value = 0.0010512900000000001
$ws.Range("H54").Value = 0.0070575000000000004
$ws.Range("I54").Value = 0.021521499999999999
$ws.Range("J54").Value = 0.067564700000000005
$ws.Range("K54").Value = 0.65173899999999996
$ws.Range("L54").Value = 0.89059500000000003
$ws.Range("M54").Value = 28.8428
# row 55
$ws.Range("D55").Value = 0.00014637499999999999
$ws.Range("E55").Value = 0.00013479099999999999
$ws.Range("F55").Value = 0.00058333300000000001
$ws.Range("G55").Value = 0.00065824999999999998
$ws.Range("H55").Value = 0.0057080799999999999
$ws.Range("I55").Value = 0.015738700000000001
$ws.Range("J55").Value = 0.10440199999999999
$ws.Range("K55").Value = 0.20430300000000001
$ws.Range("L55").Value = 1.9609300000000001
$ws.Range("M55").Value = 20.4863
# row 56
$ws.Range("D56").Value = 0.0017279999999999999
$ws.Range("E56").Value = 0.00058012500000000002
$ws.Range("F56").Value = 0.00254967
$ws.Range("G56").Value = 0.0014589200000000001
$ws.Range("H56").Value = 0.016496899999999998
$ws.Range("I56").Value = 0.0154769
$ws.Range("J56").Value = 0.13977700000000001
$ws.Range("K56").Value = 0.14511199999999999
$ws.Range("L56").Value = 2.4588000000000001
$ws.Range("M56").Value = 11.2874
# row 57
$ws.Range("D57").Value = 0.00093004100000000005
$ws.Range("E57").Value = 0.00052608400000000001
$ws.Range("F57").Value = 0.00215525
$ws.Range("G57").Value = 0.0021247100000000001
$ws.Range("H57").Value = 0.014209599999999999
$ws.Range("I57").Value = 0.0196711
$ws.Range("J57").Value = 0.16389599999999999
$ws.Range("K57").Value = 0.249253
$ws.Range("L57").Value = 3.4760399999999998
$ws.Range("M57").Value = 12.391
# row 58
$ws.Range("D58").Value = 0.012568299999999999
$ws.Range("E58").Value = 0.0079770399999999995
$ws.Range("F58").Value = 0.058606499999999999
$ws.Range("G58").Value = 0.032927400000000003
$ws.Range("H58").Value = 0.29780699999999999
$ws.Range("I58").Value = 0.417904
$ws.Range("J58").Value = 2.6355900000000001
$ws.Range("K58").Value = 3.4452600000000002
$ws.Range("L58").Value = 32.835299999999997
$ws.Range("M58").Value = 95.134200000000007

# Move the selection to match the saved view (3 strategija section header)
$ws.Range("D49:M49").Select()

Write-Output "edit complete"